# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right before "2022-Q3" (so the tab
#    order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4,
#    2021-Q3, 2021-Q2) and populate it with the Q4-2022 fund holdings table.
# 2) Prepend a corresponding "2022-Q4" row to the "总计" (totals) sheet and
#    shift the existing quarter rows down by one, renumbering the index
#    column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# =====================================================================
# 1. New "2022-Q4" worksheet
# =====================================================================
$refSheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($refSheet)
$newSheet.Name = "2022-Q4"

# Match the outline defaults (summaryBelow/summaryRight) the other sheets in
# this workbook already carry in their <sheetPr><outlinePr .../></sheetPr>.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Use an existing sheet that already has the right header/index formatting
# (bordered, bold, centered style for the header row + index column) as a
# style template - Copy() carries over both value and style, which we then
# overwrite with the real values/text below.
$template = $wb.Worksheets.Item("2022-Q3")
for ($c = 2; $c -le 8; $c++) {
    $template.Cells.Item(1, $c).Copy($newSheet.Cells.Item(1, $c))
}
for ($r = 2; $r -le 15; $r++) {
    $template.Range("A2").Copy($newSheet.Cells.Item($r, 1))
}

# Header row (row 1)
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data rows (row 2..15)
# Columns: A index(n), B code(text), C name(text), D size(text),
#          E position(text), F ratio(text), G value(text), H rank(n)
$data = @(
    @(0,  "011855", "银华长荣混合",             "10.10", "74.04", "5.00", "0.5050", 8),
    @(1,  "012186", "招商品质成长混合A",         "13.00", "92.99", "3.45", "0.4485", 10),
    @(2,  "011373", "招商前沿医疗保健股票A",      "9.48",  "92.63", "3.74", "0.3546", 10),
    @(3,  "506003", "富国科创板两年定期开放混合", "13.76", "98.91", "2.30", "0.3165", 6),
    @(4,  "630010", "华商价值精选混合",           "4.37",  "87.25", "6.57", "0.2871", 1),
    @(5,  "011598", "信澳医药健康混合",           "7.67",  "93.37", "3.48", "0.2669", 9),
    @(6,  "012187", "招商品质成长混合C",          "5.57",  "92.99", "3.45", "0.1922", 10),
    @(7,  "008978", "银华长丰混合",              "1.99",  "83.18", "5.75", "0.1144", 5),
    @(8,  "001449", "华商双驱优选灵活配置混合",    "2.26",  "77.71", "4.27", "0.0965", 4),
    @(9,  "008107", "华商医药医疗行业股票",        "1.56",  "88.55", "3.93", "0.0613", 4),
    @(10, "630006", "华商产业升级混合",           "0.86",  "88.65", "6.71", "0.0577", 1),
    @(11, "011374", "招商前沿医疗保健股票C",       "1.09",  "92.63", "3.74", "0.0408", 10),
    @(12, "005117", "金信价值精选灵活配置混合A",    "0.76",  "92.96", "2.83", "0.0215", 5),
    @(13, "005118", "金信价值精选灵活配置混合C",    "0.06",  "92.96", "2.83", "0.0017", 5)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $rec = $data[$r]

    $newSheet.Cells.Item($row, 1).Value = $rec[0]         # A index (number)

    $newSheet.Cells.Item($row, 2).Value = "'" + $rec[1]   # B fund code (forced text)
    $newSheet.Cells.Item($row, 2).Style = "Normal"

    $newSheet.Cells.Item($row, 3).Value = $rec[2]         # C fund name (text)

    $newSheet.Cells.Item($row, 4).Value = "'" + $rec[3]   # D fund size (forced text)
    $newSheet.Cells.Item($row, 4).Style = "Normal"

    $newSheet.Cells.Item($row, 5).Value = "'" + $rec[4]   # E stock position % (forced text)
    $newSheet.Cells.Item($row, 5).Style = "Normal"

    $newSheet.Cells.Item($row, 6).Value = "'" + $rec[5]   # F position ratio % (forced text)
    $newSheet.Cells.Item($row, 6).Style = "Normal"

    $newSheet.Cells.Item($row, 7).Value = "'" + $rec[6]   # G market value held (forced text)
    $newSheet.Cells.Item($row, 7).Style = "Normal"

    $newSheet.Cells.Item($row, 8).Value = $rec[7]         # H position rank (number)
}

# =====================================================================
# 2. "总计" totals sheet: prepend the 2022-Q4 row, shift the rest down
# =====================================================================
$total = $wb.Worksheets.Item("总计")

# Give the new last row's index cell (A8) the same formatting as the
# existing index cells (bordered/centered style) before overwriting values.
$total.Range("A7").Copy($total.Range("A8"))

$totalRows = @(
    @(0, "2022-Q4", 14, 2.76),
    @(1, "2022-Q3", 17, 1.97),
    @(2, "2022-Q2", 6,  0.43),
    @(3, "2022-Q1", 3,  0.32),
    @(4, "2021-Q4", 2,  0.29),
    @(5, "2021-Q3", 11, 5.32),
    @(6, "2021-Q2", 5,  1.32)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $row = $i + 2
    $rec = $totalRows[$i]
    $total.Cells.Item($row, 1).Value = $rec[0]
    $total.Cells.Item($row, 2).Value = $rec[1]
    $total.Cells.Item($row, 3).Value = $rec[2]
    $total.Cells.Item($row, 4).Value = $rec[3]
}

# Keep "总计" as the active/selected sheet (inserting a new sheet makes it
# active by default - restore the original selection).
$total.Activate()

Write-Host "2022-Q4 sheet inserted and totals sheet updated"
